# Update the "艺术表演场馆艺术演出场次观众人次" sheet:
#  - remove the oldest two years (2008年, 2009年) which occupied rows 2 and 3,
#    shifting every remaining year row up by two rows
#  - append a new row for 2021年 at the bottom (now row 13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 2008年 (row 2) and 2009年 (row 3) rows; Excel shifts rows 4.. up.
$ws.Rows("2:3").Delete() | Out-Null

# The last data row is now row 12 (2020年). Append 2021年 as row 13.
$newRow = 13

# Copy the formatting (bold + border) of the adjoining year cell before writing
# the new label so the new row matches the look of the existing year column.
$ws.Cells.Item(12, 1).Copy($ws.Cells.Item($newRow, 1))

$ws.Cells.Item($newRow, 1).Value = "2021年"
$ws.Cells.Item($newRow, 2).Value = 340.3
$ws.Cells.Item($newRow, 3).Value = 706
$ws.Cells.Item($newRow, 4).Value = 7625
$ws.Cells.Item($newRow, 6).Value = 48264.8
$ws.Cells.Item($newRow, 7).Value = 71726.5
$ws.Cells.Item($newRow, 8).Value = 15782.6
$ws.Cells.Item($newRow, 9).Value = 5453.2
$ws.Cells.Item($newRow, 10).Value = 12165.4
$ws.Cells.Item($newRow, 11).Value = 20665.4
$ws.Cells.Item($newRow, 12).Value = 62841.9
$ws.Cells.Item($newRow, 14).Value = 1305.1
$ws.Cells.Item($newRow, 15).Value = 5987.3
$ws.Cells.Item($newRow, 16).Value = 11003.4
$ws.Cells.Item($newRow, 17).Value = 83507.3
$ws.Cells.Item($newRow, 18).Value = 135.9
$ws.Cells.Item($newRow, 19).Value = 2437.7

Write-Output "ok"
